$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks one "snapshot" column per scrape; a new snapshot
# (2026-02-01 07:30:06) was inserted right before the existing "nom"
# column (CV), pushing "nom" -> CW and "url_produit" -> CW -> CX.
# Column CV = 100 (A=1).
$newColIndex = 100

# Insert a new, blank column before the current CV ("nom") column.
# This shifts CV->CW (nom) and CW->CX (url_produit), matching the diff.
$ws.Columns.Item($newColIndex).Insert()

# Header for the freshly inserted column.
$ws.Cells.Item(1, $newColIndex).Value() = "2026-02-01 07:30:06"

# The new snapshot column carries forward the last known price (same
# value/type as the previous snapshot column, now shifted to CU) for
# every data row - numeric where a price exists, blank otherwise.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $prevCell = $ws.Cells.Item($r, $newColIndex - 1)
    $newCell = $ws.Cells.Item($r, $newColIndex)
    $prevVal = $prevCell.Value()
    if ($prevVal -eq $null -or $prevVal -eq "") {
        $newCell.Value() = ""
    } else {
        $newCell.Value() = $prevVal
    }
}
